$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Names_to_change")

# Activate this sheet -> becomes the selected tab (activeTab updates, tabSelected moves here)
$ws.Activate()

$ws.Range("A2").Value = "ETH-01"
$ws.Range("B2").Value = "ETH-1"
$ws.Range("A3").Value = "ETH01"
$ws.Range("B3").Value = "ETH-1"
$ws.Range("A4").Value = "ETH1"
$ws.Range("B4").Value = "ETH-1"
$ws.Range("A5").Value = "ETH_01"
$ws.Range("B5").Value = "ETH-1"
$ws.Range("A6").Value = "ETH-02"
$ws.Range("B6").Value = "ETH-2"
$ws.Range("A7").Value = "ETH02"
$ws.Range("B7").Value = "ETH-2"
$ws.Range("A8").Value = "ETH2"
$ws.Range("B8").Value = "ETH-2"
$ws.Range("A9").Value = "ETH_02"
$ws.Range("B9").Value = "ETH-2"
$ws.Range("A10").Value = "ETH-03"
$ws.Range("B10").Value = "ETH-3"
$ws.Range("A11").Value = "ETH03"
$ws.Range("B11").Value = "ETH-3"
$ws.Range("A12").Value = "ETH3"
$ws.Range("B12").Value = "ETH-3"
$ws.Range("A13").Value = "ETH_03"
$ws.Range("B13").Value = "ETH-3"
$ws.Range("A14").Value = "ETH-04"
$ws.Range("B14").Value = "ETH-4"
$ws.Range("A15").Value = "ETH04"
$ws.Range("B15").Value = "ETH-4"
$ws.Range("A16").Value = "ETH4"
$ws.Range("B16").Value = "ETH-4"
$ws.Range("A17").Value = "ETH_04"
$ws.Range("B17").Value = "ETH-4"
$ws.Range("A18").Value = "IACA-C1"
$ws.Range("B18").Value = "IAEA-C1"
$ws.Range("A19").Value = "IAEAC1"
$ws.Range("B19").Value = "IAEA-C1"
$ws.Range("A20").Value = "IAEAC2"
$ws.Range("B20").Value = "IAEA-C2"
$ws.Range("A21").Value = "MERCK STD"
$ws.Range("B21").Value = "MERCK"
$ws.Range("A22").Value = "Merck"
$ws.Range("B22").Value = "MERCK"
$ws.Range("A23").Value = "Merck STD"
$ws.Range("B23").Value = "MERCK"
$ws.Range("A24").Value = "Rodolo Dolomite"
$ws.Range("B24").Value = "RODOLO"
$ws.Range("A25").Value = "Rodolo"
$ws.Range("B25").Value = "RODOLO"
$ws.Range("A26").Value = "RODOLO DOLOMITE"
$ws.Range("B26").Value = "RODOLO"
$ws.Range("A27").Value = "Sansa Dolomite"
$ws.Range("B27").Value = "SANSA"
$ws.Range("A28").Value = "Sansa"
$ws.Range("B28").Value = "SANSA"
$ws.Range("A29").Value = "SANSA DOLOMITE"
$ws.Range("B29").Value = "SANSA"
$ws.Range("A30").Value = "TV-03"
$ws.Range("B30").Value = "TV-04"
$ws.Range("A31").Value = "EHT-01"
$ws.Range("B31").Value = "ETH-1"
$ws.Range("A32").Value = "EHT-02"
$ws.Range("B32").Value = "ETH-2"
$ws.Range("A33").Value = "EHT-03"
$ws.Range("B33").Value = "ETH-3"
$ws.Range("A34").Value = "EHT-04"
$ws.Range("B34").Value = "ETH-4"

# Update the visible selection on the sheet to match the new data range
$ws.Range("A2:B34").Select()
